$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title
$ws.Range("A1").Value = "测试采购单"

# Row 4: 白菜 (cabbage) purchase entry
$ws.Range("A4").Value = "Sun Jul 08 2018 18:58:55 GMT+0800 (GMT+08:00)"
$ws.Range("B4").Value = "白菜"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = "斤"
$ws.Range("G4").Value = "10"
$ws.Range("H4").Value = "10"
$ws.Range("I4").Value = "采购人"
$ws.Range("J4").Value = "收验货人"
$ws.Range("K4").Value = "供货人"

# Row 5: 猪肉 (pork) purchase entry
$ws.Range("A5").Value = "Sun Jul 08 2018 19:01:02 GMT+0800 (GMT+08:00)"
$ws.Range("B5").Value = "猪肉"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "斤"
$ws.Range("G5").Value = "5"
$ws.Range("H5").Value = "5"
$ws.Range("I5").Value = "采购人"
$ws.Range("J5").Value = "收验货人"
$ws.Range("K5").Value = "供货人"

# Row 6: 白菜 (cabbage) purchase entry
$ws.Range("A6").Value = "Sun Jul 08 2018 19:05:48 GMT+0800 (GMT+08:00)"
$ws.Range("B6").Value = "白菜"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "1"
$ws.Range("F6").Value = "斤"
$ws.Range("G6").Value = "10"
$ws.Range("H6").Value = "10"
$ws.Range("I6").Value = "采购人"
$ws.Range("J6").Value = "收验货人"
$ws.Range("K6").Value = "供货人"
